$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date-serial value that was bumped by one
# day (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 45181) {
        $cell.Value = 45182
    }
}
